# Macroferia Regional de Talca - Poroto granado: weekly price update.
# A new weekly observation row is inserted before the old row 177, pushing
# the existing rows 177:203 down to 178:204 and growing the used range from
# A1:R203 to A1:R204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 177 (shifts 177:203 down to 178:204).
$ws.Rows("177:177").Insert()

# Populate the newly inserted row 177 with the new weekly record.
$ws.Cells.Item(177, 1).Value  = 5
$ws.Cells.Item(177, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(177, 3).Value  = "Maule"
$ws.Cells.Item(177, 4).Value  = 44995
$ws.Cells.Item(177, 5).Value  = 7
$ws.Cells.Item(177, 6).Value  = 100112030
$ws.Cells.Item(177, 7).Value  = "Poroto granado"
$ws.Cells.Item(177, 8).Value  = "Sin especificar"
$ws.Cells.Item(177, 9).Value  = "Primera"
$ws.Cells.Item(177, 10).Value = 200
$ws.Cells.Item(177, 11).Value = 30000
$ws.Cells.Item(177, 12).Value = 30000
$ws.Cells.Item(177, 13).Value = 30000
$ws.Cells.Item(177, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(177, 15).Value = "Región del Maule"
$ws.Cells.Item(177, 16).Value = 1200
$ws.Cells.Item(177, 17).Value = 25
$ws.Cells.Item(177, 18).Value = "Hortaliza"
